# Update "想去人数" (want-to-go count, column F) for six events across the
# "展览" and "全部类型" sheets, each incremented by 1 as a result of a
# re-scrape of the source data.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 6924
$ws1.Range("F10").Value = 16480
$ws1.Range("F11").Value = 11
$ws1.Range("F12").Value = 1620
$ws1.Range("F17").Value = 11493
$ws1.Range("F21").Value = 385

# Sheet "全部类型" (all types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 6924
$ws4.Range("F11").Value = 16480
$ws4.Range("F12").Value = 11
$ws4.Range("F13").Value = 1620
$ws4.Range("F20").Value = 11493
$ws4.Range("F24").Value = 385
